$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the placeholder empty-string values in B4:B10 ---
# These cells currently hold a shared empty string with a "quotePrefix" style.
# Copy the (plain) format from column A of the same row, then clear the
# contents so the cells become genuinely empty (keeping a normal style).
foreach ($r in 4..10) {
    $src = $ws.Range("A$r")
    $dst = $ws.Range("B$r")
    $src.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
    $dst.ClearContents()
}
$excel.CutCopyMode = 0

# --- Add the new row of data (row 14), using the same format as row 13 ---
$ws.Range("A13:C13").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A14").Value = "Teste"
$ws.Range("B14").Value = "Teste"
$ws.Range("C14").Value = "C:\Users\rapha\Desktop\Input_Python\scripts_auto_service"
$ws.Range("A14").EntireRow.RowHeight = 18.75

# --- Resize the table to include the new row ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C14"))

# --- Apply the new table style ---
$tbl.TableStyle = "TableStyleLight1"

# --- Re-fit column C now that it holds longer text ---
$ws.Columns.Item(3).AutoFit()
